$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume table (rows 2-51, columns B:E).
# A new "Frax" entry was inserted at row 41, pushing the remaining rows down
# by one and dropping the former last entry ("EOS") off the bottom of the list.
# Numeric-looking price strings are written with a leading apostrophe so Excel
# keeps them as literal text (preserving formats like trailing zeros / multiple dots).

$ws.Cells.Item(2, 2).Value = "Bitcoin"
$ws.Cells.Item(2, 3).Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Cells.Item(2, 4).Value = "27.953.89"
$ws.Cells.Item(2, 5).Value = "  +0.21%  "

$ws.Cells.Item(3, 2).Value = "Ethereum"
$ws.Cells.Item(3, 3).Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Cells.Item(3, 4).Value = "1.889.68"
$ws.Cells.Item(3, 5).Value = "  +0.52%  "

$ws.Cells.Item(4, 2).Value = "TetherUSD"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Cells.Item(4, 4).Value = "'1.014"
$ws.Cells.Item(4, 5).Value = "  +1.20%  "

$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Cells.Item(5, 4).Value = "'335.87"
$ws.Cells.Item(5, 5).Value = "  +0.81%  "

$ws.Cells.Item(6, 2).Value = "USDC"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(6, 4).Value = "'1.014"
$ws.Cells.Item(6, 5).Value = "  +1.21%  "

$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(7, 4).Value = "'0.4704"
$ws.Cells.Item(7, 5).Value = "  -0.54%  "

$ws.Cells.Item(8, 2).Value = "Cardano"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(8, 4).Value = "'0.3931"
$ws.Cells.Item(8, 5).Value = "  -1.15%  "

$ws.Cells.Item(9, 2).Value = "OKB"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(9, 4).Value = "'46.50"
$ws.Cells.Item(9, 5).Value = "  -4.58%  "

$ws.Cells.Item(10, 2).Value = "Dogecoin"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10, 4).Value = "'0.07980"
$ws.Cells.Item(10, 5).Value = "  -1.02%  "

$ws.Cells.Item(11, 2).Value = "Polygon"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(11, 4).Value = "'1.016"
$ws.Cells.Item(11, 5).Value = "  -1.21%  "

$ws.Cells.Item(12, 2).Value = "Solana"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(12, 4).Value = "'21.82"
$ws.Cells.Item(12, 5).Value = "  -0.59%  "

$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.884.58"
$ws.Cells.Item(13, 5).Value = "  -0.37%  "

$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).Value = "'5.968"
$ws.Cells.Item(14, 5).Value = "  +0.11%  "

$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).Value = "'7.152"
$ws.Cells.Item(15, 5).Value = "  -0.55%  "

$ws.Cells.Item(16, 2).Value = "BinanceUSD"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(16, 4).Value = "'1.018"
$ws.Cells.Item(16, 5).Value = "  +1.33%  "

$ws.Cells.Item(17, 2).Value = "TRON"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(17, 4).Value = "'0.06754"
$ws.Cells.Item(17, 5).Value = "  +1.94%  "

$ws.Cells.Item(18, 2).Value = "Litecoin"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(18, 4).Value = "'87.62"
$ws.Cells.Item(18, 5).Value = "  +0.41%  "

$ws.Cells.Item(19, 2).Value = "ShibaInu"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).Value = "'0.00001051"
$ws.Cells.Item(19, 5).Value = "  -0.05%  "

$ws.Cells.Item(20, 2).Value = "Avalanche"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(20, 4).Value = "'17.16"
$ws.Cells.Item(20, 5).Value = "  -0.74%  "

$ws.Cells.Item(21, 2).Value = "Dai"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(21, 4).Value = "'1.014"
$ws.Cells.Item(21, 5).Value = "  +1.20%  "

$ws.Cells.Item(22, 2).Value = "WrappedBTC"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(22, 4).Value = "27.965.49"
$ws.Cells.Item(22, 5).Value = "  -0.31%  "

$ws.Cells.Item(23, 2).Value = "Uniswap"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(23, 4).Value = "'5.499"
$ws.Cells.Item(23, 5).Value = "  -0.14%  "

$ws.Cells.Item(24, 2).Value = "Cosmos"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(24, 4).Value = "'10.98"
$ws.Cells.Item(24, 5).Value = "  -0.74%  "

$ws.Cells.Item(25, 2).Value = "Toncoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(25, 4).Value = "'2.367"
$ws.Cells.Item(25, 5).Value = "  +2.64%  "

$ws.Cells.Item(26, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(26, 4).Value = "2.111.61"
$ws.Cells.Item(26, 5).Value = "  -0.31%  "

$ws.Cells.Item(27, 2).Value = "Monero"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(27, 4).Value = "'159.35"
$ws.Cells.Item(27, 5).Value = "  +1.37%  "

$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).Value = "'20.01"
$ws.Cells.Item(28, 5).Value = "  -1.06%  "

$ws.Cells.Item(29, 2).Value = "LidoDAOToken"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(29, 4).Value = "'2.104"
$ws.Cells.Item(29, 5).Value = "  -0.13%  "

$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(30, 4).Value = "'5.511"
$ws.Cells.Item(30, 5).Value = "  -2.10%  "

$ws.Cells.Item(31, 2).Value = "BitcoinCash"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(31, 4).Value = "'121.35"
$ws.Cells.Item(31, 5).Value = "  -1.11%  "

$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).Value = "'0.09571"
$ws.Cells.Item(32, 5).Value = "  +0.06%  "

$ws.Cells.Item(33, 2).Value = "ImmutableX"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(33, 4).Value = "'0.9633"
$ws.Cells.Item(33, 5).Value = "  -1.73%  "

$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(34, 4).Value = "'3.649"
$ws.Cells.Item(34, 5).Value = "  +0.36%  "

$ws.Cells.Item(35, 2).Value = "Filecoin"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(35, 4).Value = "'5.353"
$ws.Cells.Item(35, 5).Value = "  +0.48%  "

$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "'1.362"
$ws.Cells.Item(36, 5).Value = "  -6.91%  "

$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(37, 4).Value = "'0.06134"
$ws.Cells.Item(37, 5).Value = "  +0.41%  "

$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.02251"
$ws.Cells.Item(38, 5).Value = "  -0.56%  "

$ws.Cells.Item(39, 2).Value = "TrustWalletToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(39, 4).Value = "'1.219"
$ws.Cells.Item(39, 5).Value = "  -0.77%  "

$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(40, 4).Value = "'8.191"
$ws.Cells.Item(40, 5).Value = "  -0.44%  "

$ws.Cells.Item(41, 2).Value = "Frax"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(41, 4).Value = "'1.014"
$ws.Cells.Item(41, 5).Value = "  +1.22%  "

$ws.Cells.Item(42, 2).Value = "TheSandbox"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(42, 4).Value = "'0.5963"
$ws.Cells.Item(42, 5).Value = "  -1.20%  "

$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(43, 4).Value = "'0.1900"
$ws.Cells.Item(43, 5).Value = "  -0.18%  "

$ws.Cells.Item(44, 2).Value = "Aptos"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(44, 4).Value = "'10.32"
$ws.Cells.Item(44, 5).Value = "  +0.05%  "

$ws.Cells.Item(45, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(45, 4).Value = "'1.268"
$ws.Cells.Item(45, 5).Value = "  +1.69%  "

$ws.Cells.Item(46, 2).Value = "Decentraland"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(46, 4).Value = "'0.5686"
$ws.Cells.Item(46, 5).Value = "  -0.94%  "

$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).Value = "'12.27"
$ws.Cells.Item(47, 5).Value = "  -0.24%  "

$ws.Cells.Item(48, 2).Value = "PancakeSwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(48, 4).Value = "'3.404"
$ws.Cells.Item(48, 5).Value = "  -0.32%  "

$ws.Cells.Item(49, 2).Value = "NEARProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(49, 4).Value = "'1.946"
$ws.Cells.Item(49, 5).Value = "  -0.06%  "

$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "'0.06856"
$ws.Cells.Item(50, 5).Value = "  +0.40%  "

$ws.Cells.Item(51, 2).Value = "Quant"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(51, 4).Value = "'114.00"
$ws.Cells.Item(51, 5).Value = "  -0.07%  "
